# Edit the Enterprise Support Datasheet: append trailing periods to the
# two footnote sentences in the language-support / business-hours table
# on slide 4 ("Table 6", the graphicFrame with id=25).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Locate the target table shape (graphicFrame id=25, name "Table 6") robustly
# by id rather than assuming a fixed shape index.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 25) {
        $tableShape = $candidate
        break
    }
}

$tbl = $tableShape.Table

# The footnote row is the last row; its first cell spans all 4 columns.
$cell = $tbl.Cell(3, 1)
$tr = $cell.Shape.TextFrame.TextRange

# Paragraph 1: "Language support is only available in English and Japanese " -> "...Japanese."
$paras = $tr.Paragraphs()
$para1 = $paras.Item(1)
$para1.Text = "Language support is only available in English and Japanese."

# Paragraph 3: "...P2, P3, P4 cases are limited to business hours only in Japan" -> "...Japan."
$paras = $tr.Paragraphs()
$para3 = $paras.Item(3)
$para3.Text = [char]0xA0 + "1 P2, P3, P4 cases are limited to business hours only in Japan."
